$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: site4 / tache4 / 50
$ws.Range("A5").Value = 44618
$ws.Range("A5").NumberFormat = "mmmm d yyyy"
$ws.Range("B5").Value = "site4"
$ws.Range("C5").Value = "tache4"
$ws.Range("D5").Value = 50

# Row 6: site5 / tache5 / 2
$ws.Range("A6").Value = 44618
$ws.Range("A6").NumberFormat = "mmmm d yyyy"
$ws.Range("B6").Value = "site5"
$ws.Range("C6").Value = "tache5"
$ws.Range("D6").Value = 2
